{"js": "// Fix wording mistakes in the \"if statement\" documentation paragraph:\n//   - \"there is not logic\"  -> \"there is no logic\"\n//   - \"you add between\"     -> \"you need to add between\"\n//   - \"if\u201d the logic,\"      -> \"if\u201d your logic,\"\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, \"Replace\");\n    await context.sync();\n  }\n}\n\nawait replaceOnce(\"there is not logic\", \"there is no logic\");\nawait replaceOnce(\"you add between\", \"you need to add between\");\nawait replaceOnce(\"if\\u201D the logic,\", \"if\\u201D your logic,\");\n", "ps1": "# Fix wording mistakes in the \"if statement\" documentation paragraph:\n#   - \"there is not logic\"  -> \"there is no logic\"\n#   - \"you add between\"     -> \"you need to add between\"\n#   - \"if\" the logic,\"      -> \"if\" your logic,\"\n\n$d = $word.ActiveDocument\n$quote = [char]8220\n$closeQuote = [char]8221\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n}\n\nReplace-Text \"there is not logic\" \"there is no logic\"\nReplace-Text \"you add between\" \"you need to add between\"\n\n$findQuoted = $quote + \"if\" + $closeQuote + \" the logic,\"\n$replaceQuoted = $quote + \"if\" + $closeQuote + \" your logic,\"\nReplace-Text $findQuoted $replaceQuoted\n"}
